# ---------------------------------------------------------------------------
# C1--C2-and-C3-PowerPoint.pptx - commit "Fri, Jun 05, 2020 12:04:48 AM"
#
# Two things changed in the authored commit:
#
#   1. The table on the "PLENARY - COMPLETE THE MISSING GAPS" slide (slide 16)
#      was switched from the deck's custom "Table_0" style to the built-in
#      table style {1DBAB447-05B8-46BA-BCF1-B13043E114D0}.
#
#   2. The presentation's theme colour palette was swapped from the
#      "Integral" scheme over to the standard Office colour scheme (the
#      deck's notes master already carried the stock Office theme, and the
#      edit effectively promotes that palette to the slide design as well).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 --------------------------------------

$slide = $p.Slides.Item(16)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        # Table styles are read-only as a property in this object model --
        # they must be assigned via ApplyStyle.
        $shape.Table.ApplyStyle("{1DBAB447-05B8-46BA-BCF1-B13043E114D0}")
    }
}

# --- 2. Swap the theme palette to the standard Office colours ---------------

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

function Set-ThemeRGB([int]$index, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorScheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Index map: 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3
#            8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
# dk1/lt1 (black/white) are already shared between the two palettes.
Set-ThemeRGB 3  "44546A"
Set-ThemeRGB 4  "E7E6E6"
Set-ThemeRGB 5  "5B9BD5"
Set-ThemeRGB 6  "ED7D31"
Set-ThemeRGB 7  "A5A5A5"
Set-ThemeRGB 8  "FFC000"
Set-ThemeRGB 9  "4472C4"
Set-ThemeRGB 10 "70AD47"
Set-ThemeRGB 11 "0563C1"
Set-ThemeRGB 12 "954F72"
